$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 152, pushing the existing row 152 (Dream Interpretation/
# Hecate, id 1099001) and everything below it down to row 155.
$ws.Rows("152:154").Insert()

# The "quantity" columns (AB, AH) store decimal-looking values ("1.0") as plain
# text, not numbers, in this workbook. Pre-format them as Text so assigning the
# string does not get silently coerced into the number 1.
$ws.Range("AB152:AB154").NumberFormat = "@"
$ws.Range("AH152:AH154").NumberFormat = "@"

# Row 152: Bounty Offer / Korryn
$ws.Range("A152").Value = 1047001
$ws.Range("B152").Value = 'Purple'
$ws.Range("C152").Value = 'Bounty Offer'
$ws.Range("D152").Value = '懸賞依頼'
$ws.Range("E152").Value = '현상금 의뢰'
$ws.Range("F152").Value = '悬赏委托'
$ws.Range("G152").Value = '懸賞委託'
$ws.Range("H152").Value = 'According to credible leads, a crime boss long pursued by the Public Security Bureau has fled to WhiteSands. With the outlaw''s whereabouts unknown, the Bureau decides to secretly issue a huge bounty to recruit a Sinner from WhiteSands for the capture operation.'
$ws.Range("I152").Value = '信頼できる情報提供者によると、治安局が長い間追跡してきた犯罪組織のリーダーが砂の海に逃亡したという。その後の行方が掴めていないため、治安局は密かに高額な懸賞金をかけて、砂の海出身のコンビクトに協力依頼を出した。'
$ws.Range("J152").Value = '믿을만한 소식에 따르면, 치안국에서 오랫동안 추적해 온 범죄 조직의 두목이 화이트 샌드로 도주했다고 한다. 행방이 묘연한 가운데 치안국은 비밀리에 고액 현상금을 내걸고 화이트 샌드 현지의 수감자를 모집해 검거에 나서고자 한다.'
$ws.Range("K152").Value = '据可靠线报表示，治安局追查已久的犯罪团伙头目已逃往砂海，因其行踪成谜，治安局决定秘密发布高额悬赏，招募一位砂海本地的禁闭者参与抓捕。'
$ws.Range("L152").Value = '據可靠線報表示，治安局追查已久的犯罪集團頭目已逃往砂海，因其行蹤成謎，治安局決定秘密發布高額懸賞，招募一位砂海本地的禁閉者參與抓捕。'
$ws.Range("M152").Value = 'Korryn'
$ws.Range("N152").Value = 'コリン'
$ws.Range("O152").Value = '코린'
$ws.Range("P152").Value = '柯琳'
$ws.Range("Q152").Value = '柯琳'
$ws.Range("W152").Value = 'Bronze Crystal'
$ws.Range("X152").Value = '曲銅結晶'
$ws.Range("Y152").Value = '곡선형 구리 결정'
$ws.Range("Z152").Value = '曲铜晶'
$ws.Range("AA152").Value = '曲銅晶'
$ws.Range("AB152").Value = '1.0'
$ws.Range("AC152").Value = 'Bronze Concentrate'
$ws.Range("AD152").Value = '曲銅の精鉱'
$ws.Range("AE152").Value = '정교한 곡선형 구리 광석'
$ws.Range("AF152").Value = '曲铜精矿'
$ws.Range("AG152").Value = '曲銅精礦'
$ws.Range("AH152").Value = '1.0'

# Row 153: Seasoned Mediator / Mira
$ws.Range("A153").Value = 1047002
$ws.Range("B153").Value = 'Green'
$ws.Range("C153").Value = 'Seasoned Mediator'
$ws.Range("D153").Value = '調停専門家'
$ws.Range("E153").Value = '중재 전문가'
$ws.Range("F153").Value = '调停专家'
$ws.Range("G153").Value = '調停專家'
$ws.Range("H153").Value = 'A violent dispute due to sales competition has broken out between two neighboring businesses on an Eastside shopping street. An experienced mediator with business management knowledge is urgently required to prevent further property damage from the escalating situation.'
$ws.Range("I153").Value = 'ニューシティの歩道で、隣接する二つの店舗が激しく争っている。原因は店舗間の競争の激化だ。暴力による財産の損害を軽減するため、経営と紛争調停に精通した人材を緊急派遣する必要がある。'
$ws.Range("J153").Value = '신성의 한 보행자 거리에 인접해 있는 두 점포 간에 과열된 경쟁과 갈등으로 인한 싸움이 일어났다. 경영과 중재에 능한 사람을 조속히 파견해 폭행 사건으로 인한 재산 피해를 최대한 줄여야 한다.'
$ws.Range("K153").Value = '新城某步行街，两处相邻店家发生激烈争执大打出手，起因为店铺之间的竞争与冲突，需要紧急派遣一位熟悉经营与冲突调停的人员，减少因暴力事件引起的财物损毁。'
$ws.Range("L153").Value = '新城某步行街，兩處相鄰店家發生激烈爭執大打出手，起因為店舖之間的競爭與衝突，需要緊急派遣一位熟悉經營與衝突調停的人員，減少因暴力事件引起的財物損毀。'
$ws.Range("M153").Value = 'Mira'
$ws.Range("N153").Value = 'ミラ'
$ws.Range("O153").Value = '미라'
$ws.Range("P153").Value = '米拉'
$ws.Range("Q153").Value = '米拉'
$ws.Range("W153").Value = 'Bronze Concentrate'
$ws.Range("X153").Value = '曲銅の精鉱'
$ws.Range("Y153").Value = '정교한 곡선형 구리 광석'
$ws.Range("Z153").Value = '曲铜精矿'
$ws.Range("AA153").Value = '曲銅精礦'
$ws.Range("AB153").Value = '1.0'
$ws.Range("AC153").Value = 'Bronze Raw Ore'
$ws.Range("AD153").Value = '曲銅の原鉱'
$ws.Range("AE153").Value = '거친 곡선형 구리 광석'
$ws.Range("AF153").Value = '曲铜粗矿'
$ws.Range("AG153").Value = '曲銅粗礦'
$ws.Range("AH153").Value = '1.0'

# Row 154: Management Expert / Mira + Rise
$ws.Range("A154").Value = 1047003
$ws.Range("B154").Value = 'Blue'
$ws.Range("C154").Value = 'Management Expert'
$ws.Range("D154").Value = '経営管理の人材'
$ws.Range("E154").Value = '경영관리 인재'
$ws.Range("F154").Value = '管理人才'
$ws.Range("G154").Value = '管理人才'
$ws.Range("H154").Value = 'The Public Security Bureau''s new creative merchandise has hit its lowest sales record and is making severe losses. An experienced management expert is required to provide business consultation.'
$ws.Range("I154").Value = '治安局が新たに開発した文化的なオリジナルグッズの売り上げが過去最低を記録し、損失も深刻だ。管理と経営に精通した優秀な人材の指導が急務となっている。'
$ws.Range("J154").Value = '치안국에서 새로 기획한 굿즈의 판매량이 사상 최저치를 기록해 심각한 손실이 발생했다. 경영관리 능력을 갖춘 인재를 조속히 파견해 지도해야 한다.'
$ws.Range("K154").Value = '治安局新研发的文创产品销量创史低，亏损严重，急需熟悉管理与经营的优秀人才前去指导。'
$ws.Range("L154").Value = '治安局新研發的文創產品銷量創史低，虧損嚴重，急需熟悉管理與經營的優秀人才前去指導。'
$ws.Range("M154").Value = 'Mira'
$ws.Range("N154").Value = 'ミラ'
$ws.Range("O154").Value = '미라'
$ws.Range("P154").Value = '米拉'
$ws.Range("Q154").Value = '米拉'
$ws.Range("R154").Value = 'Rise'
$ws.Range("S154").Value = 'ライズ'
$ws.Range("T154").Value = '라이즈'
$ws.Range("U154").Value = '瑞思'
$ws.Range("V154").Value = '瑞思'
$ws.Range("W154").Value = 'Bronze Concentrate'
$ws.Range("X154").Value = '曲銅の精鉱'
$ws.Range("Y154").Value = '정교한 곡선형 구리 광석'
$ws.Range("Z154").Value = '曲铜精矿'
$ws.Range("AA154").Value = '曲銅精礦'
$ws.Range("AB154").Value = '1.0'
$ws.Range("AC154").Value = 'Bronze Concentrate'
$ws.Range("AD154").Value = '曲銅の精鉱'
$ws.Range("AE154").Value = '정교한 곡선형 구리 광석'
$ws.Range("AF154").Value = '曲铜精矿'
$ws.Range("AG154").Value = '曲銅精礦'
$ws.Range("AH154").Value = '1.0'

